# Rewrite the NBA_Clinch_Dates sheet as a "Team" / "Date Eliminated" report.
# Adds a header row and re-populates every row with the updated team/date data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Team", "Date Eliminated"),
    @("Memphis Grizzlies", "Playoff"),
    @("Dallas Mavericks", "03/27/2017"),
    @("Milwaukee Bucks", "Playoff"),
    @("Los Angeles Lakers", "03/12/2017"),
    @("New York Knicks", "03/22/2017"),
    @("Phoenix Suns", "03/12/2017"),
    @("New Orleans Pelicans", "03/30/2017"),
    @("Oklahoma City Thunder", "Playoff"),
    @("Houston Rockets", "Playoff"),
    @("Boston Celtics", "Playoff"),
    @("Philadelphia 76ers", "03/22/2017"),
    @("Charlotte Hornets", "04/06/2017"),
    @("Miami Heat", "04/12/2017"),
    @("Denver Nuggets", "04/08/2017"),
    @("San Antonio Spurs", "Playoff"),
    @("Brooklyn Nets", "02/15/2017"),
    @("Minnesota Timberwolves", "03/24/2017"),
    @("Atlanta Hawks", "Playoff"),
    @("Cleveland Cavaliers", "Playoff"),
    @("Indiana Pacers", "Playoff"),
    @("Orlando Magic", "03/16/2017"),
    @("LA Clippers", "Playoff"),
    @("Sacramento Kings", "03/19/2017"),
    @("Portland Trail Blazers", "Playoff"),
    @("Golden State Warriors", "Playoff"),
    @("Chicago Bulls", "Playoff"),
    @("Utah Jazz", "Playoff"),
    @("Washington Wizards", "Playoff"),
    @("Toronto Raptors", "Playoff"),
    @("Detroit Pistons", "04/05/2017")
)

# Make sure the newly-added row (31) uses the same "text" number format as
# the rest of column B so that date-like strings ("mm/dd/yyyy") are stored
# as literal text instead of being auto-converted to date serial numbers.
$ws.Range("B31").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $team = $data[$i][0]
    $dateElim = $data[$i][1]
    $ws.Cells.Item($row, 1).Value = $team
    $ws.Cells.Item($row, 2).Value = $dateElim
}

# Header row's "Date Eliminated" cell should not carry the text number format
# that the rest of column B uses - reset it back to the default "Normal" style.
$ws.Range("B1").Style = "Normal"
